$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.5381987298234439
$ws.Range("G2").Value = 0.7114261479768182
$ws.Range("H2").Value = 0.5
$ws.Range("J2").Value = "[0 1 1 1 1 0 1 0 0 1 1 1 0 0 1 0 0 1 1 0 1 0 1 1]"

$ws.Range("D3").Value = 0.5083524027459955
$ws.Range("G3").Value = 0.6799158748411306
$ws.Range("H3").Value = 0.5833333333333334
$ws.Range("J3").Value = "[1 1 1 0 1 1 0 0 0 0 1 1 1 1 1 1 0 0 1 1 0 0 0 1]"

$ws.Range("D4").Value = 0.511332037878814

$ws.Range("D5").Value = 0.5421589494097233

$ws.Range("D6").Value = 0.5538242065587602
$ws.Range("G6").Value = 0.733242697647592
$ws.Range("J6").Value = "[0 0 1 1 1 0 0 0 1 1 1 0 0 1 1 1 0 1 0 1 1 0 1 1]"
